$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.124.02"
$ws.Range("E2").Value = "  -3.80%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.445.74"
$ws.Range("E3").Value = "  -3.96%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "525.34"
$ws.Range("E5").Value = "  -3.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.51"
$ws.Range("E6").Value = "  -8.86%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.69%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.551"
$ws.Range("E8").Value = "  -3.57%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.445.08"
$ws.Range("E9").Value = "  -5.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0987"
$ws.Range("E10").Value = "  -3.24%  "
$ws.Range("E11").Value = "  -0.56%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.28"
$ws.Range("E12").Value = "  -4.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.341"
$ws.Range("E13").Value = "  -6.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.879.68"
$ws.Range("E14").Value = "  -3.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "58.033.27"
$ws.Range("E15").Value = "  -3.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.40"
$ws.Range("E16").Value = "  -8.30%  "
$ws.Range("E17").Value = "  -4.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.453.24"
$ws.Range("E18").Value = "  -4.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.57"
$ws.Range("E19").Value = "  -6.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "317.79"
$ws.Range("E20").Value = "  -3.38%  "
$ws.Range("E21").Value = "  -5.56%  "
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.64"
$ws.Range("E23").Value = "  -5.70%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.82"
$ws.Range("E24").Value = "  -1.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.402"
$ws.Range("E25").Value = "  -8.90%  "
$ws.Range("E26").Value = "  -2.95%  "
$ws.Range("E27").Value = "  -0.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.45"
$ws.Range("E28").Value = "  -7.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0746"
$ws.Range("E29").Value = "  -8.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.44"
$ws.Range("E30").Value = "  -10.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.73"
$ws.Range("E31").Value = "  -4.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "162.68"
$ws.Range("E32").Value = "  -1.15%  "
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("E34").Value = "  -12.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.05"
$ws.Range("E35").Value = "  -4.28%  "
$ws.Range("E36").Value = "  -9.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.97"
$ws.Range("E37").Value = "  -11.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.52"
$ws.Range("E38").Value = "  -7.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.25"
$ws.Range("E39").Value = "  -2.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.49"
$ws.Range("E40").Value = "  -7.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.773"
$ws.Range("E41").Value = "  -8.42%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.997"
$ws.Range("E42").Value = "  +0.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "270.64"
$ws.Range("E43").Value = "  -11.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.96"
$ws.Range("E44").Value = "  -13.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.83"
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.582"
$ws.Range("E46").Value = "  -5.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0919"
$ws.Range("E47").Value = "  -2.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "119.50"
$ws.Range("E48").Value = "  -5.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0499"
$ws.Range("E49").Value = "  -5.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.83"
$ws.Range("E51").Value = "  -8.75%  "
